$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.212962962962963
$ws.Range("C2").Value = 0.5123456790123457
$ws.Range("J2").Value = 0.02777777777777778
$ws.Range("P2").Value = 0.1358024691358025
$ws.Range("S2").Value = 0.1111111111111111
$ws.Range("B3").Value = 0.01176470588235294
$ws.Range("C3").Value = 0.02352941176470588
$ws.Range("J3").Value = 0.05294117647058823
$ws.Range("P3").Value = 0.6882352941176471
$ws.Range("S3").Value = 0.2235294117647059
$ws.Range("O4").Value = 0.03225806451612903
$ws.Range("P4").Value = 0.6129032258064516
$ws.Range("S4").Value = 0.3548387096774194
$ws.Range("B6").Value = 0.06122448979591837
$ws.Range("D6").Value = 0.00510204081632653
$ws.Range("E6").Value = 0.00510204081632653
$ws.Range("F6").Value = 0.05612244897959184
$ws.Range("J6").Value = 0.2244897959183673
$ws.Range("O6").Value = 0.04591836734693878
$ws.Range("Q6").Value = 0.1989795918367347
$ws.Range("R6").Value = 0.08163265306122448
$ws.Range("S6").Value = 0.3214285714285715
$ws.Range("B7").Value = 0.1415525114155251
$ws.Range("D7").Value = 0.0136986301369863
$ws.Range("F7").Value = 0.045662100456621
$ws.Range("J7").Value = 0.1415525114155251
$ws.Range("O7").Value = 0.0091324200913242
$ws.Range("Q7").Value = 0.1643835616438356
$ws.Range("R7").Value = 0.0776255707762557
$ws.Range("S7").Value = 0.4063926940639269
$ws.Range("B8").Value = 0.0951276102088167
$ws.Range("D8").Value = 0.01392111368909513
$ws.Range("E8").Value = 0.002320185614849188
$ws.Range("F8").Value = 0.0580046403712297
$ws.Range("J8").Value = 0.1322505800464037
$ws.Range("O8").Value = 0.01392111368909513
$ws.Range("Q8").Value = 0.1508120649651972
$ws.Range("R8").Value = 0.1090487238979118
$ws.Range("S8").Value = 0.4245939675174014
$ws.Range("B9").Value = 0.09523809523809523
$ws.Range("D9").Value = 0.006802721088435374
$ws.Range("F9").Value = 0.04761904761904762
$ws.Range("J9").Value = 0.108843537414966
$ws.Range("O9").Value = 0.0272108843537415
$ws.Range("Q9").Value = 0.2108843537414966
$ws.Range("R9").Value = 0.1224489795918367
$ws.Range("S9").Value = 0.3809523809523809
$ws.Range("B10").Value = 0.1229641693811075
$ws.Range("D10").Value = 0.01791530944625407
$ws.Range("E10").Value = 0.0008143322475570033
$ws.Range("F10").Value = 0.07084690553745929
$ws.Range("J10").Value = 0.1343648208469055
$ws.Range("O10").Value = 0.01465798045602606
$ws.Range("Q10").Value = 0.1767100977198697
$ws.Range("R10").Value = 0.08306188925081433
$ws.Range("S10").Value = 0.3786644951140065
$ws.Range("G11").Value = 0.1428571428571428
$ws.Range("J11").Value = 0.1
$ws.Range("K11").Value = 0.2085714285714286
$ws.Range("L11").Value = 0.5428571428571428
$ws.Range("S11").Value = 0.005714285714285714
$ws.Range("G12").Value = 0.7254901960784313
$ws.Range("J12").Value = 0.1911764705882353
$ws.Range("L12").Value = 0.04411764705882353
$ws.Range("S12").Value = 0.0392156862745098
$ws.Range("G13").Value = 0.7878787878787878
$ws.Range("J13").Value = 0.1818181818181818
$ws.Range("S13").Value = 0.0303030303030303
$ws.Range("G14").Value = 0.3333333333333333
$ws.Range("J14").Value = 0.6666666666666666
$ws.Range("F15").Value = 0.02714932126696833
$ws.Range("H15").Value = 0.1900452488687783
$ws.Range("I15").Value = 0.06334841628959276
$ws.Range("J15").Value = 0.3710407239819005
$ws.Range("K15").Value = 0.06787330316742081
$ws.Range("M15").Value = 0.004524886877828055
$ws.Range("N15").Value = 0.004524886877828055
$ws.Range("O15").Value = 0.05429864253393665
$ws.Range("S15").Value = 0.2171945701357466
$ws.Range("F16").Value = 0.005649717514124294
$ws.Range("H16").Value = 0.1525423728813559
$ws.Range("I16").Value = 0.06779661016949153
$ws.Range("J16").Value = 0.4406779661016949
$ws.Range("K16").Value = 0.0903954802259887
$ws.Range("M16").Value = 0.01694915254237288
$ws.Range("N16").Value = 0.005649717514124294
$ws.Range("O16").Value = 0.0903954802259887
$ws.Range("S16").Value = 0.1299435028248588
$ws.Range("F17").Value = 0.005221932114882507
$ws.Range("H17").Value = 0.1984334203655352
$ws.Range("I17").Value = 0.04960835509138381
$ws.Range("J17").Value = 0.4412532637075718
$ws.Range("K17").Value = 0.1174934725848564
$ws.Range("M17").Value = 0.01827676240208877
$ws.Range("O17").Value = 0.08093994778067885
$ws.Range("S17").Value = 0.08877284595300261
$ws.Range("F18").Value = 0.015
$ws.Range("H18").Value = 0.18
$ws.Range("I18").Value = 0.09
$ws.Range("J18").Value = 0.425
$ws.Range("K18").Value = 0.11
$ws.Range("M18").Value = 0.025
$ws.Range("O18").Value = 0.05
$ws.Range("S18").Value = 0.105
$ws.Range("F19").Value = 0.01611535199321459
$ws.Range("H19").Value = 0.2145886344359627
$ws.Range("I19").Value = 0.07124681933842239
$ws.Range("J19").Value = 0.354537743850721
$ws.Range("K19").Value = 0.1467345207803223
$ws.Range("M19").Value = 0.01696352841391009
$ws.Range("N19").Value = 0.001696352841391009
$ws.Range("O19").Value = 0.0729431721798134
$ws.Range("S19").Value = 0.1051738761662426
